$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 17,25
$data[0,0] = '2024-06-03'
$data[0,1] = '2024-06-04'
$data[0,2] = '2024-06-19'
$data[0,3] = '한국'
$data[0,4] = '한국제14호스팩'
$data[0,5] = 4000000
$data[0,6] = 4000000
$data[0,7] = 0
$data[0,8] = 2000
$data[0,9] = 2000
$data[0,10] = 4810000
$data[0,11] = 0
$data[0,12] = 2000
$data[0,13] = '1302.62:1'
$data[0,14] = '-'
$data[0,15] = 0
$data[0,16] = 0
$data[0,17] = 0
$data[0,18] = 0
$data[0,19] = 0
$data[0,20] = 0
$data[0,21] = 0
$data[0,22] = 0
$data[0,23] = 0
$data[0,24] = '기업인수합병'
$data[1,0] = '2024-06-03'
$data[1,1] = '2024-06-04'
$data[1,2] = '2024-06-19'
$data[1,3] = '미래'
$data[1,4] = '미래에셋비전스팩5호'
$data[1,5] = 4750000
$data[1,6] = 4750000
$data[1,7] = 0
$data[1,8] = 2000
$data[1,9] = 2000
$data[1,10] = 5480000
$data[1,11] = 0
$data[1,12] = 2000
$data[1,13] = '1238.5 : 1'
$data[1,14] = '-'
$data[1,15] = 0
$data[1,16] = 0
$data[1,17] = 0
$data[1,18] = 0
$data[1,19] = 0
$data[1,20] = 0
$data[1,21] = 0
$data[1,22] = 0
$data[1,23] = 0
$data[1,24] = '기업인수합병'
$data[2,0] = '2024-05-27'
$data[2,1] = '2024-05-31'
$data[2,2] = '2024-06-19'
$data[2,3] = '한국'
$data[2,4] = '씨어스테크놀로지'
$data[2,5] = 1300000
$data[2,6] = 1300000
$data[2,7] = 0
$data[2,8] = 10500
$data[2,9] = 14000
$data[2,10] = 12293880
$data[2,11] = 0
$data[2,12] = 17000
$data[2,13] = '1084.39:1'
$data[2,14] = '2.06%'
$data[2,15] = 1386123525
$data[2,16] = 1152944128
$data[2,17] = 1885010467
$data[2,18] = -4388824706
$data[2,19] = -7990287699
$data[2,20] = -9803411085
$data[2,21] = -4430563401
$data[2,22] = -7988275689
$data[2,23] = -9916946238
$data[2,24] = '심전도검사솔루션 입원환자모니터링솔루션'
$data[3,0] = '2024-05-27'
$data[3,1] = '2024-05-31'
$data[3,2] = '2024-06-17'
$data[3,3] = '대신'
$data[3,4] = '라메디텍'
$data[3,5] = 1298000
$data[3,6] = 1298000
$data[3,7] = 0
$data[3,8] = 10400
$data[3,9] = 12700
$data[3,10] = 8650735
$data[3,11] = 0
$data[3,12] = 16000
$data[3,13] = '1115.44:1'
$data[3,14] = '9.93%'
$data[3,15] = 2055746777
$data[3,16] = 2918221978
$data[3,17] = 979078233
$data[3,18] = -3343774083
$data[3,19] = -3525649863
$data[3,20] = -1713494359
$data[3,21] = -4430074915
$data[3,22] = -8304699942
$data[3,23] = -1627684107
$data[3,24] = '초소형 레이저 의료기기 및 미용기기'
$data[4,0] = '2024-05-28'
$data[4,1] = '2024-05-29'
$data[4,2] = '2024-06-18'
$data[4,3] = 'DB'
$data[4,4] = '디비금융스팩12호'
$data[4,5] = 5000000
$data[4,6] = 5000000
$data[4,7] = 0
$data[4,8] = 2000
$data[4,9] = 2000
$data[4,10] = 5840000
$data[4,11] = 0
$data[4,12] = 2000
$data[4,13] = '1141.40:1'
$data[4,14] = '-'
$data[4,15] = 0
$data[4,16] = 0
$data[4,17] = 0
$data[4,18] = 0
$data[4,19] = 0
$data[4,20] = 0
$data[4,21] = 0
$data[4,22] = 0
$data[4,23] = 0
$data[4,24] = '기업인수 및 합병'
$data[5,0] = '2024-05-23'
$data[5,1] = '2024-05-29'
$data[5,2] = '2024-06-14'
$data[5,3] = '삼성'
$data[5,4] = '그리드위즈'
$data[5,5] = 1400000
$data[5,6] = 1400000
$data[5,7] = 0
$data[5,8] = 34000
$data[5,9] = 40000
$data[5,10] = 7942750
$data[5,11] = 0
$data[5,12] = 40000
$data[5,13] = '124.60:1'
$data[5,14] = '0.95%'
$data[5,15] = 0
$data[5,16] = 0
$data[5,17] = 0
$data[5,18] = 0
$data[5,19] = 0
$data[5,20] = 0
$data[5,21] = 0
$data[5,22] = 0
$data[5,23] = 0
$data[5,24] = '수요관리 서비스, 전기차 충전기 모뎀 등'
$data[6,0] = '2024-05-13'
$data[6,1] = '2024-05-14'
$data[6,2] = '2024-05-29'
$data[6,3] = '미래'
$data[6,4] = '미래에셋비전스팩4호'
$data[6,5] = 6650000
$data[6,6] = 6650000
$data[6,7] = 0
$data[6,8] = 2000
$data[6,9] = 2000
$data[6,10] = 8100000
$data[6,11] = 0
$data[6,12] = 2000
$data[6,13] = '1011.2:1'
$data[6,14] = '-'
$data[6,15] = 0
$data[6,16] = 0
$data[6,17] = 0
$data[6,18] = 0
$data[6,19] = 0
$data[6,20] = 0
$data[6,21] = 0
$data[6,22] = 0
$data[6,23] = 0
$data[6,24] = '기업인수목적회사(기타금융서비스)'
$data[7,0] = '2024-04-30'
$data[7,1] = '2024-05-08'
$data[7,2] = '2024-05-23'
$data[7,3] = '삼성'
$data[7,4] = '노브랜드'
$data[7,5] = 1200000
$data[7,6] = 1200000
$data[7,7] = 0
$data[7,8] = 8700
$data[7,9] = 11000
$data[7,10] = 7651263
$data[7,11] = 0
$data[7,12] = 14000
$data[7,13] = '1075.61:1'
$data[7,14] = '4.51%'
$data[7,15] = 468321534076
$data[7,16] = 555936831337
$data[7,17] = 359249623614
$data[7,18] = 22403886436
$data[7,19] = 33386727728
$data[7,20] = 10411712773
$data[7,21] = 10859975142
$data[7,22] = 29346086803
$data[7,23] = 4820429371
$data[7,24] = 'Knit, Woven 의류'
$data[8,0] = '2024-04-29'
$data[8,1] = '2024-04-30'
$data[8,2] = '2024-05-17'
$data[8,3] = 'KB'
$data[8,4] = 'KB제28호스팩'
$data[8,5] = 5000000
$data[8,6] = 5000000
$data[8,7] = 0
$data[8,8] = 2000
$data[8,9] = 2000
$data[8,10] = 5505000
$data[8,11] = 0
$data[8,12] = 2000
$data[8,13] = '1118.39:1'
$data[8,14] = '-'
$data[8,15] = 0
$data[8,16] = 0
$data[8,17] = 0
$data[8,18] = 0
$data[8,19] = 0
$data[8,20] = 0
$data[8,21] = 0
$data[8,22] = 0
$data[8,23] = 0
$data[8,24] = '기업인수합병'
$data[9,0] = '2024-04-24'
$data[9,1] = '2024-05-30'
$data[9,2] = '2024-05-17'
$data[9,3] = 'NH'
$data[9,4] = '아이씨티케이'
$data[9,5] = 1970000
$data[9,6] = 1970000
$data[9,7] = 0
$data[9,8] = 13000
$data[9,9] = 16000
$data[9,10] = 13124496
$data[9,11] = 0
$data[9,12] = 20000
$data[9,13] = '783.2:1'
$data[9,14] = '6.54%'
$data[9,15] = 0
$data[9,16] = 0
$data[9,17] = 0
$data[9,18] = 0
$data[9,19] = 0
$data[9,20] = 0
$data[9,21] = 0
$data[9,22] = 0
$data[9,23] = 0
$data[9,24] = 'PUF반도체,보안솔루션(보안반도체,정보통신모듈기기,정보통신용반도체) 제조,개발'
$data[10,0] = '2024-04-15'
$data[10,1] = '2024-04-19'
$data[10,2] = '2024-05-07'
$data[10,3] = '한국'
$data[10,4] = '코칩'
$data[10,5] = 1500000
$data[10,6] = 1500000
$data[10,7] = 0
$data[10,8] = 11000
$data[10,9] = 14000
$data[10,10] = 8503460
$data[10,11] = 0
$data[10,12] = 18000
$data[10,13] = '988.32:1'
$data[10,14] = '13.19%'
$data[10,15] = 47284698907
$data[10,16] = 38750429966
$data[10,17] = 25900014771
$data[10,18] = 7595091433
$data[10,19] = 5807002440
$data[10,20] = 3668321605
$data[10,21] = 5701880294
$data[10,22] = 4780312126
$data[10,23] = 4195570793
$data[10,24] = '소형 및 초소형 슈퍼커패시터'
$data[11,0] = '2024-04-17'
$data[11,1] = '2024-04-18'
$data[11,2] = '2024-05-07'
$data[11,3] = 'SK'
$data[11,4] = 'SK증권제12호스팩'
$data[11,5] = 3000000
$data[11,6] = 3000000
$data[11,7] = 0
$data[11,8] = 2000
$data[11,9] = 2000
$data[11,10] = 3310000
$data[11,11] = 0
$data[11,12] = 2000
$data[11,13] = '1,189.41:1'
$data[11,14] = '-'
$data[11,15] = 0
$data[11,16] = 0
$data[11,17] = 0
$data[11,18] = 0
$data[11,19] = 0
$data[11,20] = 0
$data[11,21] = 0
$data[11,22] = 0
$data[11,23] = 0
$data[11,24] = '기업인수목적 주식회사'
$data[12,0] = '2024-04-12'
$data[12,1] = '2024-04-18'
$data[12,2] = '2024-05-03'
$data[12,3] = 'KB'
$data[12,4] = '민테크'
$data[12,5] = 3000000
$data[12,6] = 3000000
$data[12,7] = 0
$data[12,8] = 6500
$data[12,9] = 8500
$data[12,10] = 21945300
$data[12,11] = 0
$data[12,12] = 10500
$data[12,13] = '946.72:1'
$data[12,14] = '4.23%'
$data[12,15] = 9576212189
$data[12,16] = 11914994171
$data[12,17] = 7285537916
$data[12,18] = 1172310325
$data[12,19] = -2762203259
$data[12,20] = -4737405164
$data[12,21] = -7460336546
$data[12,22] = -7104430732
$data[12,23] = -7501425172
$data[12,24] = '배터리 진단시스템, 배터리 시스템, 충방전 검사장비'
$data[13,0] = '2024-04-12'
$data[13,1] = '2024-04-18'
$data[13,2] = '2024-05-02'
$data[13,3] = '한국'
$data[13,4] = '디앤디파마텍'
$data[13,5] = 1100000
$data[13,6] = 1100000
$data[13,7] = 0
$data[13,8] = 22000
$data[13,9] = 26000
$data[13,10] = 10429232
$data[13,11] = 0
$data[13,12] = 33000
$data[13,13] = '848.50:1'
$data[13,14] = '10.96%'
$data[13,15] = -75676750274
$data[13,16] = -68652978862
$data[13,17] = -9506668082
$data[13,18] = -69862474811
$data[13,19] = -137025491259
$data[13,20] = 3014576074
$data[13,21] = 0
$data[13,22] = 0
$data[13,23] = 0
$data[13,24] = '대사성질환 치료제 등'
$data[14,0] = '2024-04-15'
$data[14,1] = '2024-04-16'
$data[14,2] = '2024-05-02'
$data[14,3] = '유안타'
$data[14,4] = '유안타제16호스팩'
$data[14,5] = 5150000
$data[14,6] = 5150000
$data[14,7] = 0
$data[14,8] = 2000
$data[14,9] = 2000
$data[14,10] = 5510000
$data[14,11] = 0
$data[14,12] = 2000
$data[14,13] = '1,050.42:1'
$data[14,14] = '-'
$data[14,15] = 0
$data[14,16] = 0
$data[14,17] = 0
$data[14,18] = 0
$data[14,19] = 0
$data[14,20] = 0
$data[14,21] = 0
$data[14,22] = 0
$data[14,23] = 0
$data[14,24] = '금융 지원 서비스(기업인수목적회사)'
$data[15,0] = '2024-04-08'
$data[15,1] = '2024-04-09'
$data[15,2] = '2024-04-24'
$data[15,3] = '하나'
$data[15,4] = '하나33호스팩'
$data[15,5] = 3500000
$data[15,6] = 3500000
$data[15,7] = 0
$data[15,8] = 2000
$data[15,9] = 2000
$data[15,10] = 3700000
$data[15,11] = 0
$data[15,12] = 2000
$data[15,13] = '1277.22:1'
$data[15,14] = '-'
$data[15,15] = 0
$data[15,16] = 0
$data[15,17] = 0
$data[15,18] = 0
$data[15,19] = 0
$data[15,20] = 0
$data[15,21] = 0
$data[15,22] = 0
$data[15,23] = 0
$data[15,24] = '기업인수합병'
$data[16,0] = '2024-04-04'
$data[16,1] = '2024-04-05'
$data[16,2] = '2024-04-22'
$data[16,3] = '신한'
$data[16,4] = '신한제13호스팩'
$data[16,5] = 3000000
$data[16,6] = 3000000
$data[16,7] = 0
$data[16,8] = 2000
$data[16,9] = 2000
$data[16,10] = 3620000
$data[16,11] = 0
$data[16,12] = 2000
$data[16,13] = '1337.88:1'
$data[16,14] = '-'
$data[16,15] = 0
$data[16,16] = 0
$data[16,17] = 0
$data[16,18] = 0
$data[16,19] = 0
$data[16,20] = 0
$data[16,21] = 0
$data[16,22] = 0
$data[16,23] = 0
$data[16,24] = '기타금융서비스(기업합병)'

$ws.Range("A2:Y18").Value = $data
